$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 223.88889
$ws.Range("I4").Value = 239.28572
$ws.Range("K4").Value = 239.28572
$ws.Range("M4").Value = -125.28572
$ws.Range("H32").Value = 1674.875
$ws.Range("I32").Value = 400
$ws.Range("K32").Value = 400
$ws.Range("M32").Value = -74
$ws.Range("H39").Value = 8020.154
$ws.Range("I39").Value = 373.66666
$ws.Range("J39").Value = 25224.75
$ws.Range("K39").Value = 1120.99998
$ws.Range("L39").Value = 75674.25
$ws.Range("M39").Value = -824.9999800000001
$ws.Range("N39").Value = -76266.25
$ws.Range("H80").Value = 836.3333
$ws.Range("I80").Value = 668
$ws.Range("K80").Value = 2004
$ws.Range("M80").Value = -1006
$ws.Range("H83").Value = 836.3333
$ws.Range("I83").Value = 668
$ws.Range("K83").Value = 6012
$ws.Range("M83").Value = -1020
$ws.Range("H112").Value = 5930.087
$ws.Range("J112").Value = 3682
$ws.Range("L112").Value = 11046
$ws.Range("N112").Value = -13262
$ws.Range("H137").Value = 6685.978
$ws.Range("I137").Value = 2323.2896
$ws.Range("K137").Value = 6969.8688
$ws.Range("M137").Value = -4419.8688
$ws.Range("H138").Value = 2552
$ws.Range("I138").Value = 1480.4
$ws.Range("J138").Value = 2980.64
$ws.Range("K138").Value = 4441.200000000001
$ws.Range("L138").Value = 8941.92
$ws.Range("M138").Value = 698.7999999999993
$ws.Range("N138").Value = -19221.92

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1859.8334
$ws.Range("I2").Value = 1720.3334
$ws.Range("J2").Value = 1999.3334
$ws.Range("K2").Value = 1720.3334
$ws.Range("L2").Value = 1999.3334
$ws.Range("M2").Value = -1607.3334
$ws.Range("N2").Value = -2225.3334
$ws.Range("H4").Value = 302.375
$ws.Range("I4").Value = 245.57143
$ws.Range("K4").Value = 245.57143
$ws.Range("M4").Value = -129.57143
$ws.Range("H106").Value = 91664.664
$ws.Range("J106").Value = 91664.664
$ws.Range("L106").Value = 91664.664
$ws.Range("N106").Value = -94188.664
$ws.Range("H110").Value = 5027
$ws.Range("I110").Value = 5375.3076
$ws.Range("J110").Value = 499
$ws.Range("K110").Value = 5375.3076
$ws.Range("L110").Value = 499
$ws.Range("M110").Value = -3330.3076
$ws.Range("N110").Value = -4589
$ws.Range("H116").Value = 1859.8334
$ws.Range("I116").Value = 1720.3334
$ws.Range("J116").Value = 1999.3334
$ws.Range("K116").Value = 1720.3334
$ws.Range("L116").Value = 1999.3334
$ws.Range("M116").Value = 573.6666
$ws.Range("N116").Value = -6587.3334
$ws.Range("H122").Value = 2249.5
$ws.Range("I122").Value = 2099
$ws.Range("K122").Value = 6297
$ws.Range("M122").Value = -3847
$ws.Range("H132").Value = 3683145.8
$ws.Range("I132").Value = 1133.5161
$ws.Range("K132").Value = 3400.5483
$ws.Range("M132").Value = -870.5483000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1859.8334
$ws.Range("I3").Value = 1720.3334
$ws.Range("J3").Value = 1999.3334
$ws.Range("K3").Value = 1720.3334
$ws.Range("L3").Value = 1999.3334
$ws.Range("M3").Value = -1606.3334
$ws.Range("N3").Value = -2227.3334
$ws.Range("H37").Value = 3596.7
$ws.Range("J37").Value = 4705.857
$ws.Range("L37").Value = 4705.857
$ws.Range("N37").Value = -4979.857
$ws.Range("H45").Value = 30059
$ws.Range("I45").Value = 30059
$ws.Range("K45").Value = 30059
$ws.Range("M45").Value = -29251
$ws.Range("H86").Value = 1067
$ws.Range("I86").Value = 1024.6666
$ws.Range("J86").Value = 1194
$ws.Range("K86").Value = 1024.6666
$ws.Range("L86").Value = 1194
$ws.Range("M86").Value = 98.33339999999998
$ws.Range("N86").Value = -3440
$ws.Range("H89").Value = 1067
$ws.Range("I89").Value = 1024.6666
$ws.Range("J89").Value = 1194
$ws.Range("K89").Value = 5123.333000000001
$ws.Range("L89").Value = 5970
$ws.Range("M89").Value = 492.6669999999995
$ws.Range("N89").Value = -17202
$ws.Range("H94").Value = 1293.6216
$ws.Range("I94").Value = 1284.3214
$ws.Range("K94").Value = 1284.3214
$ws.Range("M94").Value = -833.3214
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H134").Value = 52115.07
$ws.Range("J134").Value = 24729.727
$ws.Range("L134").Value = 74189.181
$ws.Range("N134").Value = -79259.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1293.2858
$ws.Range("I22").Value = 342.2857
$ws.Range("K22").Value = 342.2857
$ws.Range("M22").Value = 7.71429999999998
$ws.Range("H31").Value = 13614.962
$ws.Range("I31").Value = 1469.1818
$ws.Range("K31").Value = 1469.1818
$ws.Range("M31").Value = -1174.1818
$ws.Range("H34").Value = 13614.962
$ws.Range("I34").Value = 1469.1818
$ws.Range("K34").Value = 1469.1818
$ws.Range("M34").Value = -1267.1818
$ws.Range("H54").Value = 21121.5
$ws.Range("J54").Value = 14828.667
$ws.Range("L54").Value = 14828.667
$ws.Range("N54").Value = -16144.667
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("H63").Value = 15271
$ws.Range("J63").Value = 15271
$ws.Range("L63").Value = 15271
$ws.Range("N63").Value = -16643
$ws.Range("H66").Value = 15271
$ws.Range("J66").Value = 15271
$ws.Range("L66").Value = 45813
$ws.Range("N66").Value = -52677
$ws.Range("H122").Value = 3268.7334
$ws.Range("I122").Value = 2168.5
$ws.Range("K122").Value = 6505.5
$ws.Range("M122").Value = -4055.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1494.4546
$ws.Range("I44").Value = 491.2857
$ws.Range("K44").Value = 1473.8571
$ws.Range("M44").Value = -1075.8571
$ws.Range("H131").Value = 1477.33
$ws.Range("I131").Value = 1420.3334
$ws.Range("J131").Value = 1479.0928
$ws.Range("K131").Value = 4261.0002
$ws.Range("L131").Value = 4437.278399999999
$ws.Range("M131").Value = 778.9997999999996
$ws.Range("N131").Value = -14517.2784

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 39999.332
$ws.Range("J34").Value = 39999.332
$ws.Range("L34").Value = 39999.332
$ws.Range("N34").Value = -40535.332
$ws.Range("H52").Value = 20348.902
$ws.Range("I52").Value = 18333.334
$ws.Range("J52").Value = 20694.428
$ws.Range("K52").Value = 18333.334
$ws.Range("L52").Value = 20694.428
$ws.Range("M52").Value = -18074.334
$ws.Range("N52").Value = -21212.428
$ws.Range("H59").Value = 37300
$ws.Range("I59").Value = 35000
$ws.Range("J59").Value = 39600
$ws.Range("K59").Value = 35000
$ws.Range("L59").Value = 39600
$ws.Range("M59").Value = -34417
$ws.Range("N59").Value = -40766
$ws.Range("H68").Value = 888888
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 888888
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = 888888
$ws.Range("N68").Value = -890510
$ws.Range("L68").ClearContents()
$ws.Range("H69").Value = 70000
$ws.Range("J69").Value = 70000
$ws.Range("L69").Value = 70000
$ws.Range("N69").Value = -71498
$ws.Range("H71").Value = 888888
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 888888
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = 2666664
$ws.Range("N71").Value = -2674776
$ws.Range("L71").ClearContents()
$ws.Range("H72").Value = 70000
$ws.Range("J72").Value = 70000
$ws.Range("L72").Value = 210000
$ws.Range("N72").Value = -217488
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("H76").Value = 39999.332
$ws.Range("J76").Value = 39999.332
$ws.Range("L76").Value = 39999.332
$ws.Range("N76").Value = -40629.332
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("H79").Value = 39999.332
$ws.Range("J79").Value = 39999.332
$ws.Range("L79").Value = 39999.332
$ws.Range("N79").Value = -42183.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2999.5
$ws.Range("I22").Value = 2999
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 2999
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -2704
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 2999.5
$ws.Range("I27").Value = 2999
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 2999
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -2892
$ws.Range("N27").Value = -3214
$ws.Range("H40").Value = 2354.92
$ws.Range("I40").Value = 1585.1364
$ws.Range("K40").Value = 1585.1364
$ws.Range("M40").Value = -1449.1364
$ws.Range("H100").Value = 3316.6667
$ws.Range("I100").Value = 2955.5557
$ws.Range("J100").Value = 4400
$ws.Range("K100").Value = 2955.5557
$ws.Range("L100").Value = 4400
$ws.Range("M100").Value = -2414.5557
$ws.Range("N100").Value = -5482
$ws.Range("H105").Value = 26990
$ws.Range("J105").Value = 26990
$ws.Range("L105").Value = 26990
$ws.Range("N105").Value = -33978

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2439.5854
$ws.Range("I126").Value = 2392.0967
$ws.Range("K126").Value = 7176.2901
$ws.Range("M126").Value = -4706.2901
